$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.35395348072052
$ws.Range("B1").Value = 3.514934539794922
$ws.Range("C1").Value = 3.306158065795898
$ws.Range("D1").Value = 1.558985114097595
$ws.Range("E1").Value = 1.236286520957947
